$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new row 15 data: "101. Symmetric Tree or not"
$ws.Range("B15").Value = 101
$ws.Range("C15").Value = "Symmetric Tree"
$ws.Range("D15").Value = "Java/Python"
$ws.Range("E15").Value = "Easy"

# Copy styles from row 14 equivalents so formatting matches
$ws.Range("A15").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats (matches A15's style, same as B14)

$ws.Range("D14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("E14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# Update the active selection to E15, matching the diff
$ws.Range("E15").Select() | Out-Null
